# Commit: "Add information about topic_1 and topic_2"
#
# Slide 3 (sldId=258) of the deck has two small label text boxes that read
# "Topic 1" / "Topic 2" and two semi-transparent highlight rectangles behind
# the "subcribe ..." / "get prediction ..." callouts. The author renamed the
# labels to the snake_case form used elsewhere in the deck and nudged the two
# highlight rectangles a little so they line up better with the callouts they
# highlight.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- rename the two topic labels --------------------------------------------
$topic1 = $s.Shapes.Item("ZoneTexte 38")   # shape id 39, currently "Topic 1"
$topic1.TextFrame.TextRange.Text = "topic_1"

$topic2 = $s.Shapes.Item("ZoneTexte 39")   # shape id 40, currently "Topic 2"
$topic2.TextFrame.TextRange.Text = "topic_2"

# --- nudge the two highlight rectangles -------------------------------------
# Shape.Left / Shape.Top are exposed as single-precision points in the COM
# object model, so the literals below are the closest point values that
# round-trip to the exact target EMU offsets used in the canonical OOXML
# (914400 EMU/in, 12700 EMU/pt).

$rect1 = $s.Shapes.Item("Rectangle 1")     # shape id 2
# off x=3508188,y=3006165 -> x=3472153,y=3000068 (size unchanged)
$rect1.Left = 273.397874015748
$rect1.Top = 236.22582677165354

$rect2 = $s.Shapes.Item("Rectangle 2")     # shape id 3
# off x=8214558,y=3093651 -> x=8297941,y=3153353 (size unchanged)
$rect2.Left = 653.3811951023622
$rect2.Top = 248.29551181102363
